$d = $word.ActiveDocument
$d.Content.Find.Execute("Team Leader, Software Developer", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Lead Developer", 2)
